$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column A to hold the "Sl. No." values.
#    This shifts B->C, C->D, D->E, E->F, F->G, G->H and keeps merged cells /
#    column widths correctly shifted along with it.
$ws.Columns("A:A").Insert()

# 2. Fix the text of the "Power of 2 choices ..." header (in -> using).
$ws.Range("G1").Value = "Power of 2 choices using LRU"

# 3. Add the new "Sl. No." header.
$ws.Range("A2").Value = "Sl. No."

# 4. Fill in the new "Sl. No." data column for the existing rows.
$ws.Range("A4").Value = 1
$ws.Range("A5").Value = 2
$ws.Range("A6").Value = 3

# 5. Append two new data rows (row 7 and row 8).
$ws.Range("A7").Value = 4
$ws.Range("B7").Value = 5
$ws.Range("C7").Value = 45
$ws.Range("D7").Value = 25
$ws.Range("E7").Value = 34
$ws.Range("F7").Value = 0.28299999999999997
$ws.Range("G7").Value = 36
$ws.Range("H7").Value = 0.21

$ws.Range("A8").Value = 5
$ws.Range("B8").Value = 3
$ws.Range("C8").Value = 10
$ws.Range("D8").Value = 5
$ws.Range("E8").Value = 7
$ws.Range("F8").Value = 0.128
$ws.Range("G8").Value = 5
$ws.Range("H8").Value = 0.084000000000000005

# 6. Column widths for the newly introduced G/H columns.
$ws.Columns("G:G").ColumnWidth = 13.666666666666666
$ws.Columns("H:H").ColumnWidth = 13.998697916666666

# 7. Bold the whole header area (rows 1-2, but row 1 only from column B
#    onward since A1 is intentionally left blank/untouched), then layer
#    the specific alignments back on top (center for the merged LRU
#    headers, wrap text for the long "Limit to random number generator"
#    header).
$ws.Range("B1:H1").Font.Bold = $true
$ws.Range("A2:H2").Font.Bold = $true
$ws.Range("E1:H1").HorizontalAlignment = -4108
$ws.Range("D2").WrapText = $true

# 8. Selection as saved in the authored file.
$ws.Range("H10").Select()
